# Apply cryptocurrency price/volume updates scraped on Wed May  1 15:11:08 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (it already holds plain-text values
# such as "57.487.69" / "0.0000211" that Excel would otherwise coerce to
# numbers and normalize, dropping formatting like trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.240.24"
$ws.Range("E2").Value = "  -6.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.893.69"
$ws.Range("E3").Value = "  -4.03%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.14"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.31"
$ws.Range("E6").Value = "  -4.90%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.891.49"
$ws.Range("E8").Value = "  -4.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -9.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.67"
$ws.Range("E11").Value = "  -11.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000210"
$ws.Range("E13").Value = "  -6.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.58"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.378.81"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.905.98"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("E18").Value = "  +5.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.373.51"
$ws.Range("E19").Value = "  -6.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "400.98"
$ws.Range("E20").Value = "  -9.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.84"
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.81"
$ws.Range("E23").Value = "  -5.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.07"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.92"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.99"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.62"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0981"
$ws.Range("E33").Value = "  +4.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.41"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.907"
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("E36").Value = "  -12.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.87"
$ws.Range("E37").Value = "  -4.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.24"
$ws.Range("E38").Value = "  +5.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0616"
$ws.Range("E39").Value = "  -8.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0338"
$ws.Range("E41").Value = "  -6.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.618.10"
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "357.93"
$ws.Range("E43").Value = "  -5.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.82"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.228"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.106"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.75"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("E51").Value = "  -4.63%  "
